# "added data provider and a test case"
# Adds a new "EmailData" worksheet (after the existing "AddingToCart" sheet)
# containing a small user/password data table used as a data provider, and
# restores/adjusts the active-sheet & selection state on the other sheets.

$wb = $excel.ActiveWorkbook

# --- Add the new "EmailData" sheet as the last tab ------------------------
$sheetCount = $wb.Worksheets.Count
$ws = $wb.Worksheets.Add($null, $wb.Worksheets.Item($sheetCount))
$ws.Name = "EmailData"

# --- Fill in the data-provider table ---------------------------------------
# (entered data rows first, header row last - matches the shared-string
# insertion order captured in the workbook)
$ws.Range("A2").Value = "yahyaq91@gmail.com"
$ws.Range("B2").Value = "fddfasdfd"

$ws.Range("A3").Value = "yahyaq91@yahoo.com"
$ws.Range("B3").Value = '*&(^(*&%%^&$'

$ws.Range("A4").Value = "yahyaq91@live.com"

$ws.Range("A5").Value = "nufc@hotmail.com"
$ws.Range("B5").Value = "saddf7463"

$ws.Range("A6").Value = "realmadrid@live.com"
$ws.Range("B6").Value = "dsda%%^"

$ws.Range("A1").Value = "User"
$ws.Range("B1").Value = "Password"

# One password looks numeric, so it lands as a true number, left-aligned.
$ws.Range("B4").Value = 213243314
$ws.Range("B4").NumberFormat = "general"
$ws.Range("B4").HorizontalAlignment = -4131

# New sheet keeps its own cursor position away from the data.
$ws.Range("E12").Select() | Out-Null

# --- Re-home the active tab / selections on the pre-existing sheets --------
# "AddingToCart" no longer keeps the selection cursor/tab focus ...
$wb.Worksheets.Item("AddingToCart").Select() | Out-Null
$wb.Worksheets.Item("AddingToCart").Range("A2").Select() | Out-Null

# ... the workbook now opens back on "Filters" (selection sitting on A2).
$wb.Worksheets.Item("Filters").Select() | Out-Null
$wb.Worksheets.Item("Filters").Range("A2").Select() | Out-Null
